$d = $word.ActiveDocument

# Locate the paragraph whose entire visible text is exactly "teste"
# (a stray leftover run right before the _GoBack bookmark) and remove
# just that run's text, leaving the (now empty) paragraph - and the
# bookmark that follows it - untouched.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "teste") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark from the range so only the
    # run's text is cleared, keeping the paragraph (and the bookmark
    # that lives in it) intact.
    [void]$r.MoveEnd(1, -1)
    $r.Text = ""
}
